# Add new YouTube source / start-time data for the first 8 questions
# (rows 2-9, columns B=src and C=startTime), fix a couple of answer
# strings, and move the active-cell selection, per the "add boyidol
# some question" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "2GRP1rkE4O0"
$ws.Range("C2").Value = 13

$ws.Range("B3").Value = "I3dezFzsNss"
$ws.Range("C3").Value = 0

$ws.Range("B4").Value = "N8VRaGe3Cqs"
$ws.Range("C4").Value = 77

$ws.Range("B6").Value = "Le0CwBy4SaQ"
$ws.Range("C6").Value = 57
$ws.Range("E6").Value = "VIXX vixx"

$ws.Range("B5").Value = "ZAzWT8mRoR0"
$ws.Range("C5").Value = 80

$ws.Range("D7").Value = "방탄소년단"
$ws.Range("B7").Value = "gdZLi9oWNZg"
$ws.Range("C7").Value = 23
$ws.Range("E7").Value = "방탄 BTS bts"

$ws.Range("B8").Value = "sv53BwhUTC0"
$ws.Range("C8").Value = 60

$ws.Range("B9").Value = "0IpbvXVbBYA"
$ws.Range("C9").Value = 57

# move the selection/active cell as recorded in the saved sheet view
$ws.Range("C9").Select()
